# "Nos queremos pegar un tirito" -- fill in the second (right-hand) table of
# Hoja1 with the measured period data (column O), derive N (=O*10) and
# J (=N/2), let L (=J/K) recalc on its own, and label each row with the
# measured period string in column M. Also add the two footnote rows
# (33/34) below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row data: frequency row -> (period label for column M, raw O value)
# ---------------------------------------------------------------------------
$rowData = @(
    @{ Row = 4;  M = "4,2ms";  O = 0.118 }
    @{ Row = 5;  M = "2,9ms "; O = 0.158 }
    @{ Row = 6;  M = "2,2ms";  O = 0.194 }
    @{ Row = 7;  M = "1,9ms";  O = 0.21  }
    @{ Row = 8;  M = "720μs";  O = 0.31  }
    @{ Row = 9;  M = "360μs";  O = 0.356 }
    @{ Row = 10; M = "200μs";  O = 0.372 }
    @{ Row = 11; M = "140μs";  O = 0.374 }
    @{ Row = 12; M = "90μs";   O = 0.384 }
    @{ Row = 13; M = "60μs";   O = 0.392 }
    @{ Row = 14; M = "44μs";   O = 0.392 }
    @{ Row = 15; M = "40μs";   O = 0.392 }
    @{ Row = 16; M = "32μs";   O = 0.392 }
    @{ Row = 17; M = "7μs";    O = 0.392 }
    @{ Row = 18; M = "2,8μs";  O = 0.392 }
    @{ Row = 19; M = "1,6μs";  O = 0.392 }
    @{ Row = 20; M = "1μs";    O = 0.392 }
    @{ Row = 21; M = "800ns";  O = 0.392 }
    @{ Row = 22; M = "700ns";  O = 0.392 }
    @{ Row = 23; M = "500ns";  O = 0.392 }
    @{ Row = 24; M = "300ns";  O = 0.392 }
    @{ Row = 25; M = "100ns";  O = 0.392 }
    @{ Row = 26; M = "-";      O = 0.392 }
    @{ Row = 27; M = "-";      O = 0.392 }
    @{ Row = 28; M = "-";      O = 0.392 }
    @{ Row = 29; M = "-";      O = 0.392 }
    @{ Row = 30; M = "-";      O = 0.392 }
)

# Rows whose M cell currently carries the "alternate" font style (s=3) and
# needs to be copied onto the new footnote row (M33) before being reset
# back to the plain centered style (s=1) used everywhere else.
$ws.Range("M17").Copy()
$ws.Range("M33").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Fill in O (raw measurement), N (=O*10) and J (=N/2) for every row.
# Row 4 gets its own (non-shared) formula; rows 5:30 share one formula each,
# matching the way the existing E/L columns are already laid out.
# ---------------------------------------------------------------------------
$ws.Range("O4").Value = 0.118
$ws.Range("N4").Formula = "=O4*10"
$ws.Range("J4").Formula = "=N4/2"

for ($r = 5; $r -le 30; $r++) {
    $ws.Cells.Item($r, 15).Value = ($rowData | Where-Object { $_.Row -eq $r }).O  # column O
}
$ws.Range("N5:N30").Formula = "=O5*10"
$ws.Range("J5:J30").Formula = "=N5/2"

# ---------------------------------------------------------------------------
# Column M labels (period of the square wave at that frequency).
# ---------------------------------------------------------------------------
foreach ($entry in $rowData) {
    $ws.Cells.Item($entry.Row, 13).Value = $entry.M  # column M
}

# Cells that inherited the alternate font style (s=3) from the template get
# reset to the plain centered style (s=1) used by the rest of column M.
foreach ($r in @(17, 18, 21, 24, 26)) {
    $cell = $ws.Cells.Item($r, 13)
    $cell.Style = "Normal"
    $cell.HorizontalAlignment = -4108  # xlCenter
    $cell.VerticalAlignment = -4108    # xlCenter
}

# ---------------------------------------------------------------------------
# Footnote rows below the table.
# ---------------------------------------------------------------------------
$ws.Range("M33").Value = "80ns"
$ws.Range("M34").Value = "40ns"
$ws.Range("M34").Style = "Normal"
$ws.Range("M34").HorizontalAlignment = -4108
$ws.Range("M34").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Selection, to match where the user last clicked.
# ---------------------------------------------------------------------------
$ws.Range("O33").Select()
